$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink previously attached to C2 (mailto:juan@example.com)
$null = $ws.Range("C2").Hyperlinks.Delete()

# Row 1 headers
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "locacalizacion"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "id"
$ws.Range("E1").Value = "kind"

# Row 2 data
$ws.Range("A2").Value = "jorge"
$ws.Range("B2").Value = "18:13:14:12S"
$ws.Range("C2").Value = "jorge@email.es"
$ws.Range("D2").Value = "ID4"
$ws.Range("E2").Value = 1

# Drop the now-unused trailing columns (F:I) that used to hold
# Nacionalidad / DNI / NIF / pollingStation data
$ws.Range("F1:I1").ClearContents()
$ws.Range("F2:I2").ClearContents()

# Match the new selection left behind by the edit
$null = $ws.Range("A1:E2").Select()
